$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Miner bot / 工兵机器人): reword the effect text ---
$e3 = @"
回合结束时：本牌点数减1，然后将本牌移动到1个相邻槽位中。这之后，消灭本牌所在槽位中所有陷阱牌，再将主牌堆第1张陷阱牌放在本牌所在槽位中。本牌点数因此降低至0时，本牌不会死亡而是弃置。<br>
从手牌发动：将本牌放到房间区任意槽位。
"@
$ws.Range("E3").Value = $e3

# --- Row 4 (Mini factory / 迷你工厂): rank 2 -> 1, reword effect text ---
$ws.Range("C4").Value = 1
$e4 = @"
回合结束时：本牌点数减1。然后将弃牌堆第1张“机器人”牌放在本牌所在槽位中。本牌点数因此降低至0时，本牌不会死亡而是弃置。<br>
从手牌发动：将本牌放到房间区任意位置。
"@
$ws.Range("E4").Value = $e4

# --- Row 5 (Kamikaze bot / 自爆机器人): reword effect text (row grows taller) ---
$e5 = @"
回合结束时：本牌点数减1，然后将本牌移动到1个相邻槽位中。本牌点数因此降低至0时，使本牌所在槽位中所有其他牌点数减1，然后消灭本牌，并使玩家受到1伤害。<br>
从手牌发动：将本牌放到房间区任意位置。
"@
$ws.Range("E5").Value = $e5
$ws.Rows.Item(5).RowHeight = 57

# --- Row 6 (Sentinel bot / 哨戒机器人): rank 2 -> 1, reword effect text ---
$ws.Range("C6").Value = 1
$e6 = @"
回合结束时：本牌点数减1，然后将主牌堆第1张“机器人”牌放在本牌所在槽位中。本牌在备战区时，可以再将放置的“机器人”牌加入手牌。本牌点数因此降低至0时，本牌不会死亡而是弃置。<br>
从手牌发动：将本牌放到房间区任意位置。
"@
$ws.Range("E6").Value = $e6

# --- Update the saved selection / view state ---
$ws.Range("E7").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
